{"js": "// The paragraph \"Pro: I've thought this for a while, but\u2026\" was split across\n// three separate runs (identical formatting: sz 24, szCs 24, rtl 0). The\n// edit consolidates them into a single run carrying the full sentence.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst needle = \"thought this\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(needle) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Re-write the paragraph's range with its own (unchanged) text: this\n  // collapses the multiple runs that previously made up the sentence into\n  // a single run while preserving the run's existing character formatting.\n  const range = target.getRange();\n  range.insertText(target.text, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The paragraph \"Pro: I've thought this for a while, but\u2026\" was split across\n# three separate runs (identical formatting: sz 24, szCs 24, rtl 0). The\n# edit consolidates them into a single run carrying the full sentence.\n$d = $word.ActiveDocument\n\n# Locate the paragraph by a distinctive (plain-ASCII) substring so this\n# keeps working regardless of how the curly quote / ellipsis characters are\n# transported into this script.\n$c = $d.Content\n$c.Find.ClearFormatting()\n$found = $c.Find.Execute(\"thought this\")\n\nif ($found) {\n    $para = $c.Paragraphs(1).Range\n    # Exclude the trailing paragraph mark from the replace range.\n    $para.MoveEnd(1, -1) | Out-Null\n    $fullText = $para.Text\n\n    # Re-writing the paragraph's range with its own (unchanged) text\n    # collapses the multiple runs that previously made up the sentence\n    # into a single run, while preserving the run's character formatting.\n    $para.Find.ClearFormatting()\n    $para.Find.Replacement.ClearFormatting()\n    $para.Find.Execute($fullText, $false, $false, $false, $false, $false, $true, 1, $false, $fullText, 2) | Out-Null\n}\n"}
